$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update rows 12-14 (days 62/63 columns F, I, K)
foreach ($r in 12..14) {
    $ws.Cells.Item($r, 6).Value = "yes"   # column F
    $ws.Cells.Item($r, 9).Value = 4       # column I
    $ws.Cells.Item($r, 11).Value = 4      # column K
}

# Update view/selection: scroll back to top-left A1 and select K12
$ws.Range("K12").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
